# Apply the "update about images and image chapter" edit:
#   - Swap the labels of the two buffer boxes on the diagram:
#       "Back Buffer"    -> "Display Buffer" (typed as "Display " + "Buffer")
#       "Display Buffer" -> "Back Buffer"
#   - Refresh the cached "today" date text inside the datetimeFigureOut
#     fields on the slide master and every slide layout (19/03/2020 -> 03/08/2020).

$p = $ppt.ActivePresentation

# --- 1. Update the two rounded-rectangle labels on slide 1 -----------------

$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if (-not $shape.HasTextFrame) { continue }

    $tr = $shape.TextFrame.TextRange

    if ($tr.Text -eq "Back Buffer") {
        # Replace the leading "Back " with "Display " so the text becomes
        # "Display Buffer", split as two runs ("Display " + "Buffer") just
        # like PowerPoint does when you type over a text selection.
        $prefix = $tr.Characters(1, 5)
        $prefix.Text = "Display "
    }
    elseif ($tr.Text -eq "Display Buffer") {
        $tr.Text = "Back Buffer"
    }
}

# --- 2. Refresh the cached date text of the datetimeFigureOut fields -------

$oldDate = "19/03/2020"
$newDate = "03/08/2020"

function Update-Datefield($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if (-not $sh.HasTextFrame) { continue }
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes

for ($k = 1; $k -le $master.CustomLayouts.Count; $k++) {
    Update-DateField $master.CustomLayouts.Item($k).Shapes
}
